# Adds a trailing bold/red/Times-New-Roman space run after the final
# "comentário" run in the last paragraph of the code-sample text box on
# slide 3 ("print("Hello, World!"); # isto é um comentário").

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetParagraph = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($k = 1; $k -le $paraCount; $k++) {
            $para = $tr.Paragraphs($k)
            if ($para.Text -eq 'print("Hello, World!"); # isto é um comentário') {
                $targetSlide = $slide
                $targetShape = $shape
                $targetParagraph = $para
            }
        }
    }
}

# Appending text right after the paragraph's current end inherits the
# formatting (bold, red, Times New Roman) of the last run in that
# paragraph ("comentário"), which matches the new run's desired look.
$newRun = $targetParagraph.InsertAfter(" ")
